$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: Page changed from "Profile?" to "Profile" (done first so the now-unused
# "Profile?" shared string is freed up / the new strings below get appended in the
# same order Excel would naturally assign them)
$ws.Range("D27").Value = "Profile"

# Row 3: Page changed from "Search" to "Search & Info Display"
$ws.Range("D3").Value = "Search & Info Display"

# Row 19: New comment added
$ws.Range("E19").Value = "Not sure if this means that user is updating his/er own exercise or the general exercise database"

# Row 25: Page changed from "Search" to "profile & Info Display"
$ws.Range("D25").Value = "profile & Info Display"

# Row 5: Description text got a period added after "etc" in "(vegan, keto, etc.)"
$ws.Range("C5").Value = "When I go to view a food record in the database, I want to see the serving size, calories, dietary restriction (allergies, diabetes, etc), diet preference (vegan, keto, etc.), and whether it's a ""favorite"""

# Row 25: New comment added
$ws.Range("E25").Value = "I was imagine we will have a page for user to go directly to for all his information visualized with additional info displays about lifestyle recommendations and etc. "

# Row 11: Page changed from "Search" to "Search & Info Display"
$ws.Range("D11").Value = "Search & Info Display"

# Row 12: Page changed from "Search" to "Search & Info Display"
$ws.Range("D12").Value = "Search & Info Display"

# Row 15: Page changed from "Search" to "Search & Info Display"
$ws.Range("D15").Value = "Search & Info Display"

# Row 16: Page changed from "Search" to "Profile"
$ws.Range("D16").Value = "Profile"

# Row 17: Page changed from "Search" to "Profile"
$ws.Range("D17").Value = "Profile"

# Row 24: Page changed from "Search" to "Search & Info Display"
$ws.Range("D24").Value = "Search & Info Display"

# Row 26: Page changed from "Search" to "Search & Info Display"
$ws.Range("D26").Value = "Search & Info Display"

$ws.Range("E25").Select()
